# "fix error with dates"
#
# Two job-history dates on the CV were wrong:
#   - the current "Tech Lead" role actually started Feb 2019 (was shown as Jan 2019)
#   - consequently the previous "Lead Core-Tech Programmer" role ended Jan 2019
#     instead of running to "present"
#
# This also relocates the hidden "_GoBack" bookmark (Word always keeps exactly
# one, marking the last edit position) from the end of the "Tech Lead" line to
# the spot inside "Senior Core-Tech Programmer" where editing finished.

$d = $word.ActiveDocument

# --- 1. "Tech Lead - Splash Damage / Jan 2019 - present" -> "... / Feb 2019 - present"
$pTechLead = $d.Paragraphs.Item(14).Range
$pTechLead.Find.Execute("Jan 2019", $true, $false, $false, $false, $false, $true, 1, $false, "Feb 2019", 2)

# --- 2. "Lead Core-Tech Programmer ... / Sept 2017 - present" -> "... / Sept 2017 - Jan 2019"
$pLeadCoreTech = $d.Paragraphs.Item(18).Range
$pLeadCoreTech.Find.Execute("Sept 2017 - present", $true, $false, $false, $false, $false, $true, 1, $false, "Sept 2017 " + [char]0x2013 + " Jan 2019", 2)

# --- 3. Move the hidden _GoBack bookmark to where the edits finished, inside
#        "Senior Core-Tech Programmer" (between "Sen" and "ior").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$pSenior = $d.Paragraphs.Item(21).Range
$goBackPos = $pSenior.Start + 3
$d.Bookmarks.Add("_GoBack", $d.Range($goBackPos, $goBackPos))
